$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers for new columns I (I0) and J (IF) ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting from the existing header cell (H1) onto the new
# header cells so they share the same style index (bold, bordered,
# centered) instead of creating a brand new style.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Data values for columns I (I0) and J (IF) ---
$values = @(
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(9, 9),
    @(5, 5),
    @(7, 7),
    @(8, 8),
    @(5, 6),
    @(5, 5),
    @(8, 8),
    @(6, 7),
    @(9, 9),
    @(5, 5),
    @(8, 8),
    @(8, 8),
    @(9, 9)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
